# "added user management code"
# Refreshes the Environment table: row 2 is re-pointed from the old
# Stage/pujapowar sample row to a Dev/cbtestingteam row, and the new
# user-management columns (email id / firstname / lastname / role) that
# already existed on row 2 are backfilled onto rows 9-11 as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mirror the on-screen action of selecting the row being edited.
$ws.Rows.Item(2).Select()

# ---- Row 2: refresh environment/credential columns -----------------
$ws.Range("A2").Value = "Dev"
$ws.Range("B2").Value = "admin-lyca@yopmail.com"
$ws.Range("C2").Value = "tuZ6wG7Ysi"
$ws.Range("D2").Value = "Romania"

# A2/C2/D2 pick up the same wrap+vertical-center look already used by
# the equivalent cells on row 11 (apply per-cell - multi-area ranges
# only take the formatting on their first area in this host).
foreach ($addr in @("A2", "C2", "D2")) {
    $r = $ws.Range($addr)
    $r.VerticalAlignment = -4108
    $r.WrapText = $true
}

# ---- New user-management values (row 2 + rows 9-11) -----------------
$ws.Range("E2").Value = "cbtestingteam@yopmail.com"
$ws.Range("F2").Value = "cb"
$ws.Range("G2").Value = "testing"
$ws.Range("H2").Value = "Super Admin"

$ws.Range("E9").Value = "cbtestingteam@yopmail.com"
$ws.Range("F9").Value = "cb"
$ws.Range("G9").Value = "testing"
$ws.Range("H9").Value = "Super Admin"

$ws.Range("E10").Value = "cbtestingteam@yopmail.com"
$ws.Range("F10").Value = "cb"
$ws.Range("G10").Value = "testing"
$ws.Range("H10").Value = "Super Admin"

$ws.Range("E11").Value = "cbtestingteam@yopmail.com"
$ws.Range("F11").Value = "cb"
$ws.Range("G11").Value = "testing"
$ws.Range("H11").Value = "Super Admin"

# ---- Rebuild hyperlinks (values above are kept as display text) -----
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B11"), "mailto:admin-lyca@yopmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B10"), "mailto:admin-lyca@yopmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:admin-lyca@yopmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E9"), "mailto:cbtestingteam@yopmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E10"), "mailto:cbtestingteam@yopmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E11"), "mailto:cbtestingteam@yopmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:admin-lyca@yopmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:cbtestingteam@yopmail.com") | Out-Null

# Hyperlinks.Add() stamps a fresh font-applying style onto the cells it
# touches; flip the underline off/on to collapse each cell back onto the
# pre-existing "Hyperlink" cell style instead of leaving a duplicate one.
foreach ($addr in @("B2", "B9", "B10", "B11", "E2", "E9", "E10", "E11")) {
    $r = $ws.Range($addr)
    $r.Font.Underline = $false
    $r.Font.Underline = $true
}

# ---- New column for the "role" header --------------------------------
$ws.Columns.Item(8).ColumnWidth = 17

Write-Output "done"
